# Apply updated "dSF" (column F) values per row, per the commit:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = 9
    4  = -11
    5  = -4
    6  = -10
    7  = -5
    9  = -5
    10 = -1
    11 = -4
    12 = -4
    13 = -9
    15 = 4
    16 = 3
    17 = -2
    18 = -5
    19 = 1
    20 = 4
    21 = -3
    22 = -7
    24 = -7
    25 = -3
    27 = -4
    28 = -9
    29 = 3
    31 = 3
    32 = -2
    33 = 7
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
